# Daily attendance processing - 2026-01-11 11:32:34
#
# The "Recorded By" column (G) stores a comma-separated list of who
# recorded a session. Every cell whose value is exactly
# "System, dnasr281@gmail.com" needs to be flipped to
# "dnasr281@gmail.com, System" (list order swapped). Cells already in
# the new order, or with any other value, are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$changed = 0
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
        $changed = $changed + 1
    }
}

Write-Host "Recorded By cells swapped: $changed"
